$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh: update Price (D) and Volume(1h) (E)
# columns for each coin row, plus the two swapped rows (48/49) whose
# Coin name + Link also changed.
#
# Some Price values (e.g. "399.16", "0.0498") are valid numeric
# literals, so a plain COM .Value assignment would silently convert
# them to IEEE-754 doubles (binary round-trip noise like
# "399.16000000000003") instead of keeping them as the literal text
# the source workbook stores. For those cells we assign with a
# leading apostrophe to force text entry, then reset Style back to
# "Normal" so the quote-prefix formatting Excel applies does not
# change the cell's visible style (matches the source, which uses the
# default/unstyled format for these cells).

$ws.Range('D2').Value = '53.441.34'
$ws.Range('E2').Value = '  +3.86%  '
$ws.Range('D3').Value = '3.157.99'
$ws.Range('E3').Value = '  +3.48%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''399.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.63%  '
$ws.Range('D6').Value = '''107.98'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.72%  '
$ws.Range('D7').Value = '''0.549'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.08%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''0.611'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.41%  '
$ws.Range('E10').Value = '  +5.77%  '
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').Value = '3.648.24'
$ws.Range('E13').Value = '  +3.39%  '
$ws.Range('D14').Value = '''19.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.57%  '
$ws.Range('D15').Value = '''8.00'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.95%  '
$ws.Range('D16').Value = '''1.06'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.85%  '
$ws.Range('D17').Value = '3.149.93'
$ws.Range('E17').Value = '  +3.14%  '
$ws.Range('D18').Value = '''10.63'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').Value = '53.359.13'
$ws.Range('E19').Value = '  +3.56%  '
$ws.Range('D20').Value = '''3.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.38%  '
$ws.Range('D21').Value = '''12.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.80%  '
$ws.Range('D22').Value = '0.0₃0978'
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('D23').Value = '''71.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.25%  '
$ws.Range('D24').Value = '''271.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').Value = '''3.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.20%  '
$ws.Range('D26').Value = '''8.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').Value = '''27.79'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.25%  '
$ws.Range('D28').Value = '''7.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.37%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +2.73%  '
$ws.Range('D32').Value = '''11.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.61%  '
$ws.Range('D33').Value = '''37.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.99%  '
$ws.Range('D34').Value = '''0.0498'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.51%  '
$ws.Range('E35').Value = '  +0.56%  '
$ws.Range('D36').Value = '''50.36'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').Value = '''3.64'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.92%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').Value = '''2.77'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.88%  '
$ws.Range('D40').Value = '''4.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.32%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').Value = '''17.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.05%  '
$ws.Range('D43').Value = '''1.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('D44').Value = '''130.48'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.39%  '
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('D46').Value = '''22.45'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.65%  '
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.095.20'
$ws.Range('E48').Value = '  +3.15%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '''2.39'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').Value = '''0.0508'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +26.57%  '
$ws.Range('D51').Value = '''0.0337'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.73%  '
